$wb = $excel.ActiveWorkbook

# "Add Panels" sheet holds the cells that need updating
$ws = $wb.Worksheets.Item("Add Panels")

# Rename header labels in row 2 (F2/G2)
$ws.Range("F2").Value = "Alarm Current(A)"
$ws.Range("G2").Value = "Standby Current(A)"

# Update the panel-model cell to reflect the new model name
$ws.Range("K8").Value = "MPM800-1"

# Restore the active selection to match the saved view state
$ws.Range("B5").Select()
